# Fichas_Fuentes_Fernando_11B.docx edit
#
# 1) The "FICHAS DE REFERENCIA" title paragraph loses its red color / sz=52
#    formatting (now bold, sz=32, no color) and gains a second run
#    appending the full research-topic subtitle. A small (sz=12) empty
#    spacer paragraph is inserted right after it.
# 2) One of the redundant blank " " spacer paragraphs right after the
#    "... Seccion: B" line is removed.

$d = $word.ActiveDocument

# --- Part 1: rewrite the "FICHAS DE REFERENCIA" heading paragraph ---

$heading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "FICHAS DE REFERENCIA") {
        $heading = $para
        break
    }
}

$titleXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>FICHAS DE REFERENCIA</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>: Influencia de los videojuegos en los estudiantes de secundaria del Colegio Champagnat en el período de enero de 2020 a julio de 2022</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="12"/><w:szCs w:val="12"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

if ($heading -ne $null) {
    $heading.Range.InsertXML($titleXml)
} else {
    Write-Host "WARNING: 'FICHAS DE REFERENCIA' heading paragraph not found."
}

# --- Part 2: remove one of the redundant blank spacer paragraphs ---
# They all read as a single space run (sz 32/32) right after
# "... Seccion: B"; the diff drops exactly the first one of that run.

$nbsp = [string][char]0x00A0
$target = $null
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $trimmed = $para.Range.Text.TrimEnd([char]13)
    if ($trimmed -eq $nbsp -or $trimmed -eq " ") {
        $prev = $d.Paragraphs.Item($i - 1)
        $prevTrimmed = $prev.Range.Text.TrimEnd([char]13)
        if ($prevTrimmed.EndsWith("B")) {
            $target = $para
            break
        }
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
} else {
    Write-Host "WARNING: redundant spacer paragraph after 'Seccion: B' not found."
}
